# "Generate Report for Handoff"
#
# The workbook tracks localization status of e2e\a.md and e2e\b.md across
# an Overview sheet and one sheet per target locale (zh-cn, de-de).
# b.md has just been handed off for localization, so its status/row needs
# to reflect the newly generated handoff package (new xlf file names,
# new timestamps, new status text, and a new "out of date" error detail
# on the per-locale sheets).

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$overviewDate = "2016-08-24 18:46:20"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ff0668f20ad14cb6d9c4d3a87c7c6fcba1e8ef0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce0e57c663b854362cf952695ec4aaa12452221e/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $overviewDate

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-24 18:46:15"
$wsZhCn.Range("P3").Value = $errorDetail

# widen the Error Detail column now that it holds a long message
$wsZhCn.Range("P:P").ColumnWidth = 39.2

# ---------------------------------------------------------------------
# de-de sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = $overviewDate
$wsDeDe.Range("P3").Value = $errorDetail

# widen the Error Detail column now that it holds a long message
$wsDeDe.Range("P:P").ColumnWidth = 39.2
